# Fruta / hortaliza, semanal
# Insert a new weekly record as row 184, pushing the existing rows 184-203
# down to 185-204.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 184 - this shifts rows 184:203 down
# to 185:204 and keeps all their data/formatting intact.
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new record.
$ws.Cells.Item(184, 1).Value = 3
$ws.Cells.Item(184, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(184, 3).Value = "Coquimbo"
$ws.Cells.Item(184, 4).Value = 44449
$ws.Cells.Item(184, 5).Value = 5
$ws.Cells.Item(184, 6).Value = 100112032
$ws.Cells.Item(184, 7).Value = "Zapallo italiano"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 215
$ws.Cells.Item(184, 11).Value = 14500
$ws.Cells.Item(184, 12).Value = 15500
$ws.Cells.Item(184, 13).Value = 15023
$ws.Cells.Item(184, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(184, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(184, 16).Value = 215
$ws.Cells.Item(184, 17).Value = 70
$ws.Cells.Item(184, 18).Value = "Hortaliza"
